$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "0.17"
$ws.Range("C2").Value = "44.29***"
$ws.Range("D2").Value = "-0.89"

$ws.Range("B3").Value = "-0.01"
$ws.Range("C3").Value = "2.21***"
$ws.Range("D3").Value = "0.46***"

$ws.Range("B4").Value = "-0.09"
$ws.Range("C4").Value = "0.98"
$ws.Range("D4").Value = "0.82*"
